$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header rename: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2. Column A data values (rows 2-14): generation counts -> normalized MaxFES fractions
$newA = 0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1
for ($i = 0; $i -lt $newA.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $newA[$i]
}

# 3. Recompute the "Mean" column (currently AZ holds "Run 50" data, BA holds "Mean")
#    New mean is the average of Run 0..Run 49 (columns B:AY), excluding Run 50.
$newMean = 55.76953916, 48.24416067, 11.84234314, 0.27504109, 0.23723307, 0.21212458, 0.20133319, 0.18746576, 0.17947303, 0.17445075, 0.16750882, 0.16310485, 0.15761084
for ($i = 0; $i -lt $newMean.Length; $i++) {
    $ws.Cells.Item($i + 2, 52).Value = $newMean[$i]
}

# 4. Drop the "Run 50" run entirely: delete column BA (old Mean column, now redundant)
#    and relabel column AZ's header from "Run 50" to "Mean".
$ws.Columns("BA:BA").Delete()
$ws.Range("AZ1").Value = "Mean"
